{"js": "// Apply the commit's edits:\n// 1. \"A EAP tem poucos itens, talvez tenha sido pouco subdividida.\" -> \"5\"\n// 2. Append \"1.1.1 \" before \"Projeto instrucional\" (two occurrences)\n// 3. Append \"1.1.2.1 \" before \"Componente relat\u00f3rio\" (two occurrences)\n\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"A EAP tem poucos itens, talvez tenha sido pouco subdividida.\",\n    replace: \"5\"\n  },\n  {\n    find: \"Existem pacotes de trabalho n\u00e3o detalhados no dicionario da EAP. Exemplo: Projeto instrucional\",\n    replace: \"Existem pacotes de trabalho n\u00e3o detalhados no dicionario da EAP. Exemplo: 1.1.1 Projeto instrucional\"\n  },\n  {\n    find: \"Existem pacotes de trabalho sem tamanho estimado. Exemplo: Projeto instrucional\",\n    replace: \"Existem pacotes de trabalho sem tamanho estimado. Exemplo: 1.1.1 Projeto instrucional\"\n  },\n  {\n    find: \"Existem pacotes de trabalho sem unidade de tamanho definida. Exemplo: Projeto instrucional\",\n    replace: \"Existem pacotes de trabalho sem unidade de tamanho definida. Exemplo: 1.1.1 Projeto instrucional\"\n  },\n  {\n    find: \"Existem pacotes de trabalho n\u00e3o detalhados no dicionario da EAP. Exemplo: Componente relat\u00f3rio\",\n    replace: \"Existem pacotes de trabalho n\u00e3o detalhados no dicionario da EAP. Exemplo: 1.1.2.1 Componente relat\u00f3rio\"\n  },\n  {\n    find: \"Existem pacotes de trabalho sem atividades derivadas. Exemplo: Componente relat\u00f3rio\",\n    replace: \"Existem pacotes de trabalho sem atividades derivadas. Exemplo: 1.1.2.1 Componente relat\u00f3rio\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the commit's edits by matching whole paragraphs of text and\n# rewriting their Range.Text (excluding the paragraph-mark character so\n# the existing run's rPr / xml:space formatting is preserved).\n#\n# 1. \"A EAP tem poucos itens, talvez tenha sido pouco subdividida.\" -> \"5\"\n# 2. Prefix \"1.1.1 \" onto \"Projeto instrucional\" (3 paragraphs)\n# 3. Prefix \"1.1.2.1 \" onto \"Componente relat\u00f3rio\" (2 paragraphs)\n\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"A EAP tem poucos itens, talvez tenha sido pouco subdividida.\" = \"5\";\n    \"Existem pacotes de trabalho n\u00e3o detalhados no dicionario da EAP. Exemplo: Projeto instrucional\" = \"Existem pacotes de trabalho n\u00e3o detalhados no dicionario da EAP. Exemplo: 1.1.1 Projeto instrucional\";\n    \"Existem pacotes de trabalho sem tamanho estimado. Exemplo: Projeto instrucional\" = \"Existem pacotes de trabalho sem tamanho estimado. Exemplo: 1.1.1 Projeto instrucional\";\n    \"Existem pacotes de trabalho sem unidade de tamanho definida. Exemplo: Projeto instrucional\" = \"Existem pacotes de trabalho sem unidade de tamanho definida. Exemplo: 1.1.1 Projeto instrucional\";\n    \"Existem pacotes de trabalho n\u00e3o detalhados no dicionario da EAP. Exemplo: Componente relat\u00f3rio\" = \"Existem pacotes de trabalho n\u00e3o detalhados no dicionario da EAP. Exemplo: 1.1.2.1 Componente relat\u00f3rio\";\n    \"Existem pacotes de trabalho sem atividades derivadas. Exemplo: Componente relat\u00f3rio\" = \"Existem pacotes de trabalho sem atividades derivadas. Exemplo: 1.1.2.1 Componente relat\u00f3rio\";\n}\n\nforeach ($p in $d.Paragraphs) {\n    $rng = $p.Range\n    [void]$rng.MoveEnd(1, -1)\n    $old = $rng.Text\n    if ($replacements.ContainsKey($old)) {\n        $rng.Text = $replacements[$old]\n    }\n}\n"}
